# The deck ships two theme parts:
#   theme1.xml -> "Office Theme" (used by the Notes Master)
#   theme2.xml -> "Integral"     (used by the Slide Master / main design)
# The commit swaps the two themes' contents wholesale (Integral <-> Office
# Theme). The only functional difference between the two theme definitions
# is the 10 scheme colours (dk2, lt2, accent1-6, hlink, folHlink) - fonts and
# format scheme are identical - so recreate the swap by writing the "Office
# Theme" palette into the live (slide-master) colour scheme via the Design /
# ColorScheme object model.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$colorScheme = $design.SlideMaster.Theme.ThemeColorScheme

# Index -> (scheme slot, new RGB as a VBA-style 0x00BBGGRR integer)
# New values are the current "Office Theme" srgbClr values (hex RRGGBB
# converted to the BGR-packed integer PowerPoint's ColorFormat.RGB uses).
$colorScheme.Item(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477    # folHlink 954F72
